$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.01"
$ws.Range("E2").Value = "'1.88%"
$ws.Range("G2").Value = "'8"
$ws.Range("D3").Value = "'41.57"
$ws.Range("E3").Value = "'5.52%"
$ws.Range("G3").Value = "'8"
$ws.Range("D4").Value = "'5.628"
$ws.Range("E4").Value = "'-1.56%"
$ws.Range("G4").Value = "'8"
$ws.Range("D5").Value = "'0.08182"
$ws.Range("E5").Value = "'2.14%"
$ws.Range("G5").Value = "'8"
$ws.Range("D6").Value = "'2.028"
$ws.Range("E6").Value = "'2.16%"
$ws.Range("G6").Value = "'8"
$ws.Range("D7").Value = "'8.744"
$ws.Range("E7").Value = "'1.44%"
$ws.Range("G7").Value = "'8"
$ws.Range("D8").Value = "'4.503"
$ws.Range("E8").Value = "'-1.00%"
$ws.Range("G8").Value = "'8"
$ws.Range("D9").Value = "'2.952"
$ws.Range("E9").Value = "'0.07%"
$ws.Range("G9").Value = "'8"
$ws.Range("D10").Value = "'0.9220"
$ws.Range("E10").Value = "'-0.59%"
$ws.Range("G10").Value = "'8"
$ws.Range("D11").Value = "'0.1275"
$ws.Range("E11").Value = "'0.21%"
$ws.Range("G11").Value = "'8"
$ws.Range("D12").Value = "'0.1955"
$ws.Range("E12").Value = "'0.22%"
$ws.Range("G12").Value = "'8"
$ws.Range("D13").Value = "'0.09262"
$ws.Range("E13").Value = "'1.56%"
$ws.Range("G13").Value = "'8"
$ws.Range("D14").Value = "'0.03808"
$ws.Range("E14").Value = "'3.23%"
$ws.Range("G14").Value = "'8"
$ws.Range("E15").Value = "'1.12%"
$ws.Range("G15").Value = "'8"
$ws.Range("D16").Value = "'0.001302"
$ws.Range("E16").Value = "'1.22%"
$ws.Range("G16").Value = "'8"
$ws.Range("D17").Value = "'0.006158"
$ws.Range("E17").Value = "'-1.87%"
$ws.Range("G17").Value = "'8"
$ws.Range("G18").Value = "'8"
$ws.Range("E19").Value = "'2.76%"
$ws.Range("G19").Value = "'8"
$ws.Range("D20").Value = "'0.3484"
$ws.Range("E20").Value = "'-1.49%"
$ws.Range("G20").Value = "'8"
$ws.Range("D21").Value = "'8.289"
$ws.Range("E21").Value = "'-4.93%"
$ws.Range("G21").Value = "'8"
$ws.Range("D22").Value = "'0.1373"
$ws.Range("E22").Value = "'0.02%"
$ws.Range("G22").Value = "'8"
$ws.Range("D23").Value = "'0.2411"
$ws.Range("E23").Value = "'-1.61%"
$ws.Range("G23").Value = "'8"
$ws.Range("D24").Value = "'0.04405"
$ws.Range("E24").Value = "'-0.40%"
$ws.Range("G24").Value = "'8"
$ws.Range("D25").Value = "'0.001266"
$ws.Range("E25").Value = "'-0.05%"
$ws.Range("G25").Value = "'8"
$ws.Range("D26").Value = "'0.004298"
$ws.Range("E26").Value = "'-4.74%"
$ws.Range("G26").Value = "'8"
$ws.Range("E27").Value = "'2.53%"
$ws.Range("G27").Value = "'8"
$ws.Range("G28").Value = "'8"
$ws.Range("G29").Value = "'8"
$ws.Range("G30").Value = "'8"
$ws.Range("G31").Value = "'8"
$ws.Range("G32").Value = "'8"
$ws.Range("G33").Value = "'8"
$ws.Range("G34").Value = "'8"
$ws.Range("G35").Value = "'8"
$ws.Range("G36").Value = "'8"
$ws.Range("G37").Value = "'8"
$ws.Range("G38").Value = "'8"
$ws.Range("D39").Value = "'0.02773"
$ws.Range("E39").Value = "'11.61%"
$ws.Range("G39").Value = "'8"
$ws.Range("D40").Value = "'0.05400"
$ws.Range("E40").Value = "'1.79%"
$ws.Range("G40").Value = "'8"
$ws.Range("D41").Value = "'0.007691"
$ws.Range("E41").Value = "'3.02%"
$ws.Range("G41").Value = "'8"
$ws.Range("D42").Value = "'0.1419"
$ws.Range("E42").Value = "'1.26%"
$ws.Range("G42").Value = "'8"
$ws.Range("D43").Value = "'0.008963"
$ws.Range("E43").Value = "'-6.25%"
$ws.Range("G43").Value = "'8"
$ws.Range("D44").Value = "'0.002142"
$ws.Range("E44").Value = "'1.05%"
$ws.Range("G44").Value = "'8"
$ws.Range("D45").Value = "'0.01164"
$ws.Range("E45").Value = "'5.79%"
$ws.Range("G45").Value = "'8"
$ws.Range("D46").Value = "'0.00006629"
$ws.Range("E46").Value = "'-1.95%"
$ws.Range("G46").Value = "'8"
$ws.Range("E47").Value = "'-0.13%"
$ws.Range("G47").Value = "'8"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.003194"
$ws.Range("E48").Value = "'7.22%"
$ws.Range("G48").Value = "'8"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.002281"
$ws.Range("E49").Value = "'-0.56%"
$ws.Range("G49").Value = "'8"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("G50").Value = "'8"
$ws.Range("E51").Value = "'-0.13%"
$ws.Range("G51").Value = "'8"
